# The deck's single Slide Master carried the "Integral" theme
# (ppt/theme/theme2.xml) and the Notes Master carried "Office Theme"
# (ppt/theme/theme1.xml). The edit swaps these two color palettes, so the
# slides end up themed with the plain "Office" palette, while "Integral"'s
# colors move over to the notes-master theme part.
#
# PowerPoint's object model exposes theme colors for writing via
# Slide.ThemeColorScheme (and SlideRange.ThemeColorScheme), which targets
# the Slide Master's theme. We push the "Office Theme" color scheme onto
# it one RGB value at a time, in the standard clrScheme slot order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$officeThemeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $hexColor = $officeThemeColors[$i]
    $r = [int](($hexColor -band 0xFF0000) -shr 16)
    $g = [int](($hexColor -band 0x00FF00) -shr 8)
    $b = [int]($hexColor -band 0x0000FF)
    # VBA-style RGB() encoding used by the ThemeColor.RGB property: r + g*256 + b*65536
    $bgr = $r + ($g * 256) + ($b * 65536)
    $tcs.Colors($i + 1).RGB = $bgr
}
